# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# --- Hoja "1er Parcial" ---
$ws1 = $wb.Worksheets.Item("1er Parcial")

# Fila 5 (6AEM)
$ws1.Range("E5").Value = 24
$ws1.Range("F5").Value = 12
$ws1.Range("G5").Value = 66.67
$ws1.Range("H5").Value = 33.33
$ws1.Range("I5").Value = 8.6
$ws1.Range("J5").Value = 12
$ws1.Range("K5").Value = 33.33

# Fila 6 (6ASM)
$ws1.Range("E6").Value = 19
$ws1.Range("F6").Value = 7
$ws1.Range("G6").Value = 73.08
$ws1.Range("H6").Value = 26.92
$ws1.Range("I6").Value = 8.1
$ws1.Range("J6").Value = 7
$ws1.Range("K6").Value = 26.92

# --- Hoja "2o Parcial" ---
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Fila 5 (6AEM)
$ws2.Range("E5").Value = 19
$ws2.Range("F5").Value = 17
$ws2.Range("G5").Value = 52.78
$ws2.Range("H5").Value = 47.22
$ws2.Range("I5").Value = 8.699999999999999
$ws2.Range("J5").Value = 17
$ws2.Range("K5").Value = 47.22

# Fila 6 (6ASM)
$ws2.Range("E6").Value = 19
$ws2.Range("F6").Value = 7
$ws2.Range("G6").Value = 73.08
$ws2.Range("H6").Value = 26.92
$ws2.Range("I6").Value = 8.300000000000001
$ws2.Range("J6").Value = 7
$ws2.Range("K6").Value = 26.92

# Fila 8 (6BEV)
$ws2.Range("E8").Value = 7
$ws2.Range("F8").Value = 15
$ws2.Range("G8").Value = 31.82
$ws2.Range("H8").Value = 68.18000000000001
$ws2.Range("I8").Value = 8
$ws2.Range("J8").Value = 15
$ws2.Range("K8").Value = 68.18000000000001

# --- Hoja "3er Parcial" ---
$ws3 = $wb.Worksheets.Item("3er Parcial")

# Fila 5 (6AEM)
$ws3.Range("E5").Value = 24
$ws3.Range("F5").Value = 12
$ws3.Range("G5").Value = 66.67
$ws3.Range("H5").Value = 33.33
$ws3.Range("I5").Value = 8.699999999999999
$ws3.Range("J5").Value = 12
$ws3.Range("K5").Value = 33.33

# Fila 6 (6ASM)
$ws3.Range("E6").Value = 19
$ws3.Range("F6").Value = 7
$ws3.Range("G6").Value = 73.08
$ws3.Range("H6").Value = 26.92
$ws3.Range("I6").Value = 8.4
$ws3.Range("J6").Value = 7
$ws3.Range("K6").Value = 26.92

# Fila 8 (6BEV)
$ws3.Range("I8").Value = 7.4
